$wb = $excel.ActiveWorkbook

# GPS信息表 (2nd sheet): "主键" -> "主键，自动递增", cursor moves to G2
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("G2").Value = "主键，自动递增"
$ws2.Range("G2").Select()

# 车辆信息表（vehicle）(3rd sheet): "主键" -> "主键，自动递增", becomes the active sheet,
# cursor moves to G2 (was G8)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("G2").Value = "主键，自动递增"
$ws3.Activate()
$ws3.Range("G2").Select()
